$d = $word.ActiveDocument

# Find all paragraphs with text "[placeholder]" and handle them in order.
$placeholderParas = @()
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "[placeholder]") {
        $placeholderParas += $p
    }
}

# First placeholder -> replace text with "Snake"
$firstPara = $placeholderParas[0]
$firstRange = $firstPara.Range
$firstRange.Find.Execute("[placeholder]", $true, $false, $false, $false, $false, $true, 1, $false, "Snake", 2)

# Second placeholder paragraph: delete entirely, along with the following empty paragraph
$secondPara = $placeholderParas[1]
$secondIndex = $secondPara.Range.Start
$followingPara = $secondPara.Next()

$deleteStart = $secondPara.Range.Start
$deleteEnd = $followingPara.Range.End

$deleteRange = $d.Range($deleteStart, $deleteEnd)
$deleteRange.Delete()
